# Update "Forecast Comparison" sheet: insert a Week_Start_Date column, change
# the Week labels to the un-padded form, and mark is_holiday_week as boolean.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new, blank column before the existing column B ("ASIN"); this
# shifts ASIN, MyForecast, Amazon Mean/P70/P80/P90 Forecast, Product Title
# and is_holiday_week each one column to the right (B->C, C->D, ..., I->J).
$ws.Columns("B").Insert()

# Treat the new column as text so the date-like strings we are about to
# write aren't auto-converted into date serial numbers.
$ws.Columns("B").NumberFormat = "@"

# Header for the newly inserted column.
$ws.Range("B1").Value = "Week_Start_Date"

# Column A week labels change from zero-padded ("W01".."W16") to the
# un-padded form ("W1".."W16").
$weeks = @("W1","W2","W3","W4","W5","W6","W7","W8","W9","W10","W11","W12","W13","W14","W15","W16")
for ($i = 0; $i -lt $weeks.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $weeks[$i]
}

# Populate the new Week_Start_Date column (B) with each week's start date.
$dates = @("2025-01-05","2025-01-12","2025-01-19","2025-01-26","2025-02-02","2025-02-09","2025-02-16","2025-02-23","2025-03-02","2025-03-09","2025-03-16","2025-03-23","2025-03-30","2025-04-06","2025-04-13","2025-04-20")
for ($i = 0; $i -lt $dates.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $dates[$i]
}

# is_holiday_week (now column J after the insert) becomes a boolean column.
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 10).Value = $false
}
